$wb = $excel.ActiveWorkbook

# Worksheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 53.6
$ws.Range("I12").Value = 23
$ws.Range("J12").Value = 99.5
$ws.Range("K12").Value = 23
$ws.Range("L12").Value = 99.5
$ws.Range("M12").Value = 147
$ws.Range("N12").Value = -439.5
$ws.Range("H32").Value = 873.2857
$ws.Range("I32").Value = 979.6667
$ws.Range("J32").Value = 793.5
$ws.Range("K32").Value = 979.6667
$ws.Range("L32").Value = 793.5
$ws.Range("M32").Value = -653.6667
$ws.Range("N32").Value = -1445.5
$ws.Range("H132").Value = 3985.0303
$ws.Range("I132").Value = 3682.9822
$ws.Range("K132").Value = 11048.9466
$ws.Range("M132").Value = -8518.946599999999
$ws.Range("H137").Value = 6243.4546
$ws.Range("I137").Value = 7201.048
$ws.Range("J137").Value = 3149.6924
$ws.Range("K137").Value = 21603.144
$ws.Range("L137").Value = 9449.0772
$ws.Range("M137").Value = -19053.144
$ws.Range("N137").Value = -14549.0772
$ws.Range("H138").Value = 1838.4681
$ws.Range("I138").Value = 1275.3611
$ws.Range("J138").Value = 3681.3635
$ws.Range("K138").Value = 3826.0833
$ws.Range("L138").Value = 11044.0905
$ws.Range("M138").Value = 1313.9167
$ws.Range("N138").Value = -21324.0905
$ws.Range("H141").Value = 4030.8484
$ws.Range("I141").Value = 4155.4814
$ws.Range("K141").Value = 12466.4442
$ws.Range("M141").Value = -7286.444199999998

# Worksheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1544.83
$ws.Range("I32").Value = 1461.5275
$ws.Range("K32").Value = 1461.5275
$ws.Range("M32").Value = -1174.5275
$ws.Range("H61").Value = 13035.167
$ws.Range("I61").Value = 14622.2
$ws.Range("J61").Value = 5100
$ws.Range("K61").Value = 14622.2
$ws.Range("L61").Value = 5100
$ws.Range("M61").Value = -14410.2
$ws.Range("N61").Value = -5524
$ws.Range("H74").Value = 2880.9211
$ws.Range("I74").Value = 1124.6428
$ws.Range("K74").Value = 1124.6428
$ws.Range("M74").Value = -250.6428000000001
$ws.Range("H77").Value = 2880.9211
$ws.Range("I77").Value = 1124.6428
$ws.Range("K77").Value = 5623.214
$ws.Range("M77").Value = -1255.214
$ws.Range("H132").Value = 4392.8086
$ws.Range("I132").Value = 4370.8433
$ws.Range("J132").Value = 4458.706
$ws.Range("K132").Value = 13112.5299
$ws.Range("L132").Value = 13376.118
$ws.Range("M132").Value = -10582.5299
$ws.Range("N132").Value = -18436.118
$ws.Range("H136").Value = 13035.167
$ws.Range("I136").Value = 14622.2
$ws.Range("J136").Value = 5100
$ws.Range("K136").Value = 43866.60000000001
$ws.Range("L136").Value = 15300
$ws.Range("M136").Value = -41316.60000000001
$ws.Range("N136").Value = -20400

# Worksheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11312.333
$ws.Range("I86").Value = 12288.875
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 12288.875
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -11165.875
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 11312.333
$ws.Range("I89").Value = 12288.875
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 61444.375
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -55828.375
$ws.Range("N89").Value = -28732
$ws.Range("H105").Value = 3068.0952
$ws.Range("I105").Value = 1429.1333
$ws.Range("J105").Value = 7165.5
$ws.Range("K105").Value = 1429.1333
$ws.Range("L105").Value = 7165.5
$ws.Range("M105").Value = 317.8667
$ws.Range("N105").Value = -10659.5
$ws.Range("H134").Value = 8153.689
$ws.Range("I134").Value = 8168.683
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 24506.049
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -21971.049
$ws.Range("N134").Value = -29070

# Worksheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 85353.914
$ws.Range("I16").Value = 1718.5
$ws.Range("J16").Value = 252624.75
$ws.Range("K16").Value = 1718.5
$ws.Range("L16").Value = 252624.75
$ws.Range("M16").Value = -1431.5
$ws.Range("N16").Value = -253198.75
$ws.Range("H31").Value = 2439.5813
$ws.Range("I31").Value = 1700.2162
$ws.Range("J31").Value = 6999
$ws.Range("K31").Value = 1700.2162
$ws.Range("L31").Value = 6999
$ws.Range("M31").Value = -1405.2162
$ws.Range("N31").Value = -7589
$ws.Range("H34").Value = 2439.5813
$ws.Range("I34").Value = 1700.2162
$ws.Range("J34").Value = 6999
$ws.Range("K34").Value = 1700.2162
$ws.Range("L34").Value = 6999
$ws.Range("M34").Value = -1498.2162
$ws.Range("N34").Value = -7403
$ws.Range("H58").Value = 1241.683
$ws.Range("I58").Value = 953.7778
$ws.Range("J58").Value = 1796.9286
$ws.Range("K58").Value = 953.7778
$ws.Range("L58").Value = 1796.9286
$ws.Range("M58").Value = -750.7778
$ws.Range("N58").Value = -2202.9286
$ws.Range("H86").Value = 13295.846
$ws.Range("I86").Value = 8446.666999999999
$ws.Range("J86").Value = 14750.6
$ws.Range("K86").Value = 8446.666999999999
$ws.Range("L86").Value = 14750.6
$ws.Range("M86").Value = -7323.666999999999
$ws.Range("N86").Value = -16996.6
$ws.Range("H89").Value = 13295.846
$ws.Range("I89").Value = 8446.666999999999
$ws.Range("J89").Value = 14750.6
$ws.Range("K89").Value = 42233.335
$ws.Range("L89").Value = 73753
$ws.Range("M89").Value = -36617.335
$ws.Range("N89").Value = -84985
$ws.Range("H113").Value = 85353.914
$ws.Range("I113").Value = 1718.5
$ws.Range("J113").Value = 252624.75
$ws.Range("K113").Value = 1718.5
$ws.Range("L113").Value = 252624.75
$ws.Range("M113").Value = 451.5
$ws.Range("N113").Value = -256964.75
$ws.Range("H122").Value = 9235.25
$ws.Range("I122").Value = 16524.5
$ws.Range("J122").Value = 1946
$ws.Range("K122").Value = 49573.5
$ws.Range("L122").Value = 5838
$ws.Range("M122").Value = -47123.5
$ws.Range("N122").Value = -10738
$ws.Range("H132").Value = 16333.5
$ws.Range("I132").Value = 1956.1052
$ws.Range("J132").Value = 46685.777
$ws.Range("K132").Value = 5868.3156
$ws.Range("L132").Value = 140057.331
$ws.Range("M132").Value = -3338.3156
$ws.Range("N132").Value = -145117.331
$ws.Range("H134").Value = 2345.2
$ws.Range("I134").Value = 3176.7778
$ws.Range("J134").Value = 1664.8182
$ws.Range("K134").Value = 9530.3334
$ws.Range("L134").Value = 4994.4546
$ws.Range("M134").Value = -6995.3334
$ws.Range("N134").Value = -10064.4546
$ws.Range("H136").Value = 1241.683
$ws.Range("I136").Value = 953.7778
$ws.Range("J136").Value = 1796.9286
$ws.Range("K136").Value = 2861.3334
$ws.Range("L136").Value = 5390.7858
$ws.Range("M136").Value = -311.3334
$ws.Range("N136").Value = -10490.7858

# Worksheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 174.83333
$ws.Range("I26").Value = 90.72727
$ws.Range("J26").Value = 1100
$ws.Range("K26").Value = 272.18181
$ws.Range("L26").Value = 3300
$ws.Range("M26").Value = 15.81818999999996
$ws.Range("N26").Value = -3876
$ws.Range("H56").Value = 7931.75
$ws.Range("I56").Value = 7931.75
$ws.Range("K56").Value = 7931.75
$ws.Range("M56").Value = -7401.75
$ws.Range("H113").Value = 7138.684
$ws.Range("J113").Value = 7510.278
$ws.Range("L113").Value = 22530.834
$ws.Range("N113").Value = -26870.834
$ws.Range("H136").Value = 5517.826
$ws.Range("I136").Value = 995.17645
$ws.Range("J136").Value = 18332
$ws.Range("K136").Value = 2985.52935
$ws.Range("L136").Value = 54996
$ws.Range("M136").Value = 2114.47065
$ws.Range("N136").Value = -65196
$ws.Range("H137").Value = 7138.067
$ws.Range("J137").Value = 11324.75
$ws.Range("L137").Value = 33974.25
$ws.Range("N137").Value = -44174.25

# Worksheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 15377026
$ws.Range("I3").Value = 30000550
$ws.Range("J3").Value = 753501
$ws.Range("K3").Value = 30000550
$ws.Range("L3").Value = 753501
$ws.Range("M3").Value = -30000434
$ws.Range("N3").Value = -753733
$ws.Range("H11").Value = 11396971
$ws.Range("I11").Value = 12445343
$ws.Range("J11").Value = 3009999.8
$ws.Range("K11").Value = 12445343
$ws.Range("L11").Value = 3009999.8
$ws.Range("M11").Value = -12445204
$ws.Range("N11").Value = -3010277.8
$ws.Range("H12").Value = 5155714.5
$ws.Range("I12").Value = 3235555.5
$ws.Range("J12").Value = 16676667
$ws.Range("K12").Value = 3235555.5
$ws.Range("L12").Value = 16676667
$ws.Range("M12").Value = -3235415.5
$ws.Range("N12").Value = -16676947
$ws.Range("H14").Value = 3343555.8
$ws.Range("I14").Value = 3343555.8
$ws.Range("K14").Value = 3343555.8
$ws.Range("M14").Value = -3343387.8

# Worksheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 16108.177
$ws.Range("I40").Value = 17236.32
$ws.Range("K40").Value = 17236.32
$ws.Range("M40").Value = -17100.32
$ws.Range("H136").Value = 3869.2444
$ws.Range("I136").Value = 1712.2916
$ws.Range("J136").Value = 6334.3335
$ws.Range("K136").Value = 5136.8748
$ws.Range("L136").Value = 19003.0005
$ws.Range("M136").Value = -2586.8748
$ws.Range("N136").Value = -24103.0005

# Worksheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 20000
$ws.Range("J7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("N7").Value = -20226
$ws.Range("H122").Value = 28799.773
$ws.Range("I122").Value = 5016
$ws.Range("J122").Value = 37718.688
$ws.Range("K122").Value = 15048
$ws.Range("L122").Value = 113156.064
$ws.Range("M122").Value = -12598
$ws.Range("N122").Value = -118056.064
$ws.Range("H136").Value = 180195.97
$ws.Range("I136").Value = 208815.33
$ws.Range("J136").Value = 3709.9167
$ws.Range("K136").Value = 626445.99
$ws.Range("L136").Value = 11129.7501
$ws.Range("M136").Value = -623895.99
$ws.Range("N136").Value = -16229.7501
